$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 51
    6  = 174
    7  = 257
    11 = 38
    14 = 1311
    17 = 435
    22 = 1327
    23 = 3281
    24 = 26
    27 = 1066
    28 = 72
    29 = 1656
    31 = 456
    32 = 43
    36 = 622
    37 = 426
    38 = 23
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
